$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 6.044599999999999
$ws.Range("B10").Value = 5.163200000000002
$ws.Range("B12").Value = 5.039999999999999
$ws.Range("B18").Value = 7.277099999999995
$ws.Range("B25").Value = 5.927300000000002
$ws.Range("B37").Value = 8.760700000000002
$ws.Range("B55").Value = 6.747799999999993
$ws.Range("B68").Value = 4.665799999999997
$ws.Range("B77").Value = 9.136100000000004
$ws.Range("B78").Value = 9.483100000000002
$ws.Range("B79").Value = 8.799500000000005
$ws.Range("B80").Value = 9.3971
$ws.Range("B81").Value = 5.6363
$ws.Range("B82").Value = 5.4895
$ws.Range("B84").Value = 6.232299999999998
$ws.Range("B101").Value = 9.163199999999993
$ws.Range("B102").Value = 8.362900000000003
